$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-10-25 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-10-26 Sunday", 2)

$d.Content.Find.Execute("72×21=", $true, $false, $false, $false, $false, $true, 1, $false, "52×27=", 2)
$d.Content.Find.Execute("62×22=", $true, $false, $false, $false, $false, $true, 1, $false, "17×72=", 2)
$d.Content.Find.Execute("74×48=", $true, $false, $false, $false, $false, $true, 1, $false, "93×68=", 2)
$d.Content.Find.Execute("33×47=", $true, $false, $false, $false, $false, $true, 1, $false, "91×50=", 2)
$d.Content.Find.Execute("65×53=", $true, $false, $false, $false, $false, $true, 1, $false, "26×66=", 2)

$d.Content.Find.Execute("79×54=", $true, $false, $false, $false, $false, $true, 1, $false, "56×95=", 2)
$d.Content.Find.Execute("21×90=", $true, $false, $false, $false, $false, $true, 1, $false, "69×26=", 2)
$d.Content.Find.Execute("22×87=", $true, $false, $false, $false, $false, $true, 1, $false, "25×62=", 2)
$d.Content.Find.Execute("64×74=", $true, $false, $false, $false, $false, $true, 1, $false, "94×92=", 2)
$d.Content.Find.Execute("58×98=", $true, $false, $false, $false, $false, $true, 1, $false, "17×17=", 2)

$d.Content.Find.Execute("61×39=", $true, $false, $false, $false, $false, $true, 1, $false, "73×91=", 2)
$d.Content.Find.Execute("94×58=", $true, $false, $false, $false, $false, $true, 1, $false, "21×16=", 2)
$d.Content.Find.Execute("87×28=", $true, $false, $false, $false, $false, $true, 1, $false, "13×21=", 2)
$d.Content.Find.Execute("36×46=", $true, $false, $false, $false, $false, $true, 1, $false, "18×70=", 2)
$d.Content.Find.Execute("65×30=", $true, $false, $false, $false, $false, $true, 1, $false, "54×18=", 2)

$d.Content.Find.Execute("88×60=", $true, $false, $false, $false, $false, $true, 1, $false, "59×56=", 2)
$d.Content.Find.Execute("65×57=", $true, $false, $false, $false, $false, $true, 1, $false, "90×93=", 2)
$d.Content.Find.Execute("82×81=", $true, $false, $false, $false, $false, $true, 1, $false, "12×82=", 2)
$d.Content.Find.Execute("87×47=", $true, $false, $false, $false, $false, $true, 1, $false, "69×75=", 2)
$d.Content.Find.Execute("24×61=", $true, $false, $false, $false, $false, $true, 1, $false, "65×26=", 2)

$d.Content.Find.Execute("82×56=", $true, $false, $false, $false, $false, $true, 1, $false, "76×38=", 2)
$d.Content.Find.Execute("85×15=", $true, $false, $false, $false, $false, $true, 1, $false, "52×93=", 2)
$d.Content.Find.Execute("33×25=", $true, $false, $false, $false, $false, $true, 1, $false, "76×82=", 2)
$d.Content.Find.Execute("12×36=", $true, $false, $false, $false, $false, $true, 1, $false, "85×87=", 2)
$d.Content.Find.Execute("89×90=", $true, $false, $false, $false, $false, $true, 1, $false, "90×62=", 2)
